$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add two new review rows (35 and 36), mirroring the existing data layout ---
# Start by duplicating the formatting of the last existing row (34) down into
# the two new rows so fonts/alignment/borders/number-formats match.
$ws.Range("A34:G34").Copy()
$ws.Range("A35:G36").PasteSpecial(-4122)

# Row 35: com.hamxa.shaynachim / bitcoin review from gregneri12@gmail.com
$ws.Range("A35").Value = "com.hamxa.shaynachim"
$ws.Range("B35").Value = "bitcoin"
$ws.Range("C35").Value = "gregneri12@gmail.com"
$ws.Range("D35").Value = "halachme@gmail.com"
$ws.Range("E35").Value = "27/5/2019 15:59"
$ws.Range("F35").Value = "why 5 star? because its legit...very good app"
$ws.Range("G35").Value = "no"

# Row 36: com.hamxa.shaynachim / bitcoin review from armonravid2@gmail.com
$ws.Range("A36").Value = "com.hamxa.shaynachim"
$ws.Range("B36").Value = "bitcoin"
$ws.Range("C36").Value = "armonravid2@gmail.com"
$ws.Range("D36").Value = "armonravid@gmail.com"
$ws.Range("E36").Value = "27/5/2019 15:59"
$ws.Range("F36").Value = "I've given this app a 5 star before because it never fail to work properly "
$ws.Range("G36").Value = "no"

# --- Hyperlink the email addresses in columns C and D, same as earlier rows ---
$ws.Hyperlinks.Add($ws.Range("C35"), "mailto:gregneri12@gmail.com", [Type]::Missing, [Type]::Missing, "gregneri12@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D35"), "mailto:halachme@gmail.com", [Type]::Missing, [Type]::Missing, "halachme@gmail.com")
$ws.Hyperlinks.Add($ws.Range("C36"), "mailto:armonravid2@gmail.com", [Type]::Missing, [Type]::Missing, "armonravid2@gmail.com")
$ws.Hyperlinks.Add($ws.Range("D36"), "mailto:armonravid@gmail.com", [Type]::Missing, [Type]::Missing, "armonravid@gmail.com")

# Adding a hyperlink re-styles the cell with the built-in "Hyperlink" look;
# restore the original column formatting (from row 34) over C35:D36 so the
# cells keep matching the rest of the table.
$ws.Range("C34:D34").Copy()
$ws.Range("C35:D36").PasteSpecial(-4122)

# Review rows wrap their text and are taller than the header rows.
$ws.Rows.Item(35).RowHeight = 24
$ws.Rows.Item(36).RowHeight = 24
$ws.Range("F35").WrapText = $true
$ws.Range("F36").WrapText = $true

# Move the selection to the newly added rows, matching where the edit left off.
$ws.Range("C36:D36").Select()
